$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cited_by_count column holds text digits (e.g. "18"); keep it text-typed
# (matching the original inlineStr cells) instead of letting Excel coerce
# the numeric-looking string into a real number.

# --- Row 2: cited_by_count 18 -> 21 ---
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "21"

# --- Row 3: cited_by_count 7 -> 8 ---
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "8"

# --- Row 4: affiliation text updated; cited_by_count 5 -> 6 ---
$ws.Range("B4").Value = "Arizona State University; Johns Hopkins School of Medicine; Sleep and Circadian Research Laboratory, Department of Psychiatry, University of Michigan, Ann Arbor, MI, USA; College of Health Solutions, Arizona State University, Phoenix, USA; Division of Pulmonary, Allergy, Critical Care and Sleep Medicine, University of Arizona, Tucson, USA; Edson College of Nursing and Health Innovation, Arizona State University, Health North, Suite 301, 550 N 3rd Street, Phoenix, AZ, 85004, USA; Edson College of Nursing and Health Innovation, Arizona State University, Phoenix, USA; Edson College of Nursing and Health Innovation, Arizona State University, Phoenix, USA; Edson College of Nursing and Health Innovation, Arizona State University, Phoenix, USA; Pain Research & Intervention Center of Excellence, University of Florida, Gainesville, USA; Edson College of Nursing and Health Innovation, Arizona State University, Phoenix, USA"
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = "6"

# --- Row 8: affiliation text trimmed (duplicate center references removed) ---
$ws.Range("B8").Value = "Department of Medical Social Sciences, Northwestern University Feinberg School of Medicine, Chicago, IL, United States; Division of Community and Systems Health Science, University of Arizona College of Nursing, Tucson, AZ, United States; Department of Medical Social Sciences, Northwestern University Feinberg School of Medicine, Chicago, IL, United States; Department of Neurology, Northwestern University Feinberg School of Medicine, Chicago, IL, United States; Department of Medical Social Sciences, Northwestern University Feinberg School of Medicine, Chicago, IL, United States; Robert H. Lurie Comprehensive Cancer Center of Northwestern University, Chicago, IL, United States; Department of Psychiatry, University of California, San Diego, San Diego, CA, United States; Robert H. Lurie Comprehensive Cancer Center of Northwestern University, Chicago, IL, United States; Department of Medicine and Psychology and Sylvester Comprehensive Cancer Center, University of Miami, Miami, FL, United States; Aarhus Institute of Advanced Studies, Aarhus University, Aarhus, Denmark; Department of Medical Social Sciences, Northwestern University Feinberg School of Medicine, Chicago, IL, United States; Department of Neurology, Northwestern University Feinberg School of Medicine, Chicago, IL, United States; Division of Pulmonary, Allergy, Critical Care, and Sleep Medicine, Department of Medicine, University of Arizona, Tucson, AZ, United States; Division of Community and Systems Health Science, University of Arizona College of Nursing, Tucson, AZ, United States; Department of Medical Social Sciences, Northwestern University Feinberg School of Medicine, Chicago, IL, United States; Department of Neurology, Northwestern University Feinberg School of Medicine, Chicago, IL, United States; Nox Health, Suwanee, GA, United States"

# --- Rows 11-13: cyclic rotation of records (11<-12, 12<-13, 13<-11) ---
$cols = @("A","B","C","D","H","P")
$orig11 = @{}
$orig12 = @{}
$orig13 = @{}
foreach ($col in $cols) {
    $orig11[$col] = $ws.Range($col + "11").Value2
    $orig12[$col] = $ws.Range($col + "12").Value2
    $orig13[$col] = $ws.Range($col + "13").Value2
}
foreach ($col in $cols) {
    $ws.Range($col + "11").Value = $orig12[$col]
    $ws.Range($col + "12").Value = $orig13[$col]
    $ws.Range($col + "13").Value = $orig11[$col]
}

# --- Rows 14-15: host_organization N/A -> Springer Nature ---
$ws.Range("G14").Value = "Springer Nature"
$ws.Range("G15").Value = "Springer Nature"

# --- Row 16: author name corrections ---
$ws.Range("A16").Value = "Mladen Jergović, Makiko Watanabe, Ruchika Bhat, Christopher P Coplen, Sandip Ashok Sonar, Rachel Wong, Yvonne Castaneda, Lisa Davidson, Mrinalini Kala, Rachel C. Wilson, Homer L. Twigg, Kenneth S. Knox, Heidi E Erickson, Craig Weinkauf, Christian Bime, Billie Bixby, Sairam Parthasarathy, Jarrod Mosier, Bonnie LaFleur, Deepta Bhattacharya, Janko Nikolich‐Žugich"

Write-Output "Done applying DOM and Banner author id updates."
